# Add new daily rows (row 450 through row 479, 2021-06-01 .. 2021-06-30) to each
# of the four "province" sheets, matching the upstream data upload.
#
# For every sheet:
#   - Row 450 gets a date in column A, a new case-count in column C, and a
#     7-day rolling AVERAGE formula in column D (continuing the existing
#     AVERAGE(C(r-6):Cr) pattern already used down column D).
#   - Rows 451-479 only get a date in column A (no counts reported yet for
#     those future days).

$wb = $excel.ActiveWorkbook

# sheetName -> value to put in C450 for that sheet
$newValues = @{
    "Nuovi casi"        = 9
    "Deceduti"          = 0
    "Dimessi   Guariti" = 106
    "Ricoveri"          = 20
}

$firstNewDate = 44348   # 2021-06-01 (Excel serial date)
$lastNewRow   = 479
$firstNewRow  = 450

foreach ($sheetName in $newValues.Keys) {
    $ws = $wb.Worksheets.Item($sheetName)
    $cValue = $newValues[$sheetName]

    # --- Row 450: date + count + rolling-average formula ---------------
    $ws.Cells.Item(449, 1).Copy()
    $ws.Cells.Item($firstNewRow, 1).PasteSpecial(-4122)
    $ws.Cells.Item(449, 4).Copy()
    $ws.Cells.Item($firstNewRow, 4).PasteSpecial(-4122)

    $ws.Cells.Item($firstNewRow, 1).Value = $firstNewDate
    $ws.Cells.Item($firstNewRow, 3).Value = $cValue
    $ws.Cells.Item($firstNewRow, 4).Formula = "=AVERAGE(C444:C450)"

    # --- Rows 451-479: date only ----------------------------------------
    $dateSerial = $firstNewDate + 1
    for ($r = $firstNewRow + 1; $r -le $lastNewRow; $r++) {
        $ws.Cells.Item(449, 1).Copy()
        $ws.Cells.Item($r, 1).PasteSpecial(-4122)
        $ws.Cells.Item($r, 1).Value = $dateSerial
        $dateSerial = $dateSerial + 1
    }
}

# --- View-state tweaks (scroll position / active cell / active tab) ------

$ws1 = $wb.Worksheets.Item("Nuovi casi")
[void]$ws1.Range("A450:D450").Select()

$ws2 = $wb.Worksheets.Item("Deceduti")
[void]$ws2.Range("A450:D450").Select()

$ws3 = $wb.Worksheets.Item("Dimessi   Guariti")
[void]$ws3.Range("A450:D450").Select()

$ws4 = $wb.Worksheets.Item("Ricoveri")
[void]$ws4.Range("A450:A479").Select()

# "Dimessi   Guariti" becomes the active/visible sheet (was "Ricoveri").
[void]$ws3.Activate()
